$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate order record that used to sit in row 2
# (6444a94280ca9ea5d90965f0 / audi / london / prince homes / 2023-04-12 / 2023-04-24 / 42900.00 / 0)
# Every row below shifts up by one.
$ws.Rows.Item(2).Delete()

# Append three new order records at the bottom of the table (rows 7-9).
# Columns A-G are text in this sheet, so force text formatting before
# assigning values (avoids ids/dates/amounts being parsed as numbers/dates).
$newRows = @(
    @(7, "6444eefa0181a0655cdd09bd", "audi", "america", "mumbai", "2023-04-20", "2023-04-27", "26400.00"),
    @(8, "6444ef7c0181a0655cdd09db", "Toyota Innova ", "rajkot", "mumbai", "2023-04-14", "2023-04-28", "33000.00"),
    @(9, "644663592cb7d38a677d1201", "mini cooper", "surat", "mansi ", "2023-04-17", "2023-04-26", "55000.00")
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Range("A$row`:G$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $r[1]
    $ws.Range("B$row").Value = $r[2]
    $ws.Range("C$row").Value = $r[3]
    $ws.Range("D$row").Value = $r[4]
    $ws.Range("E$row").Value = $r[5]
    $ws.Range("F$row").Value = $r[6]
    $ws.Range("G$row").Value = $r[7]
    $ws.Range("H$row").Value = 0
}
